$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A="I'm gonna grab some coffee."; B="나 커피 좀 사러 갈 거야."},
    @{Row=3;  A="I didn't mean to hurt your feelings."; B="네 기분을 상하게 하려던 건 아니었어."},
    @{Row=4;  A="Are you sure you want to quit?"; B="정말 그만두고 싶은 거 맞아?"},
    @{Row=5;  A="I'm here to see Mr. Kim."; B="김 선생님을 뵈러 왔습니다."},
    @{Row=6;  A="It looks like it's going to rain."; B="비가 올 것 같아."},
    @{Row=7;  A="How about we go out for dinner?"; B="우리 저녁 먹으러 나가는 거 어때?"},
    @{Row=8;  A="Why don't you take a break?"; B="좀 쉬는 게 어때?"},
    @{Row=9;  A="I used to play the piano."; B="예전엔 피아노를 치곤 했어."},
    @{Row=10; A="I'd like to make a reservation."; B="예약을 하고 싶은데요."},
    @{Row=11; A="Do you mind if I sit here?"; B="여기 앉아도 될까요?"},
    @{Row=12; A="What if I fail the test?"; B="시험에 떨어지면 어떡해?"},
    @{Row=13; A="It's time to say goodbye."; B="이제 헤어질 시간이야."},
    @{Row=14; A="There is no need to worry."; B="걱정할 필요 없어."},
    @{Row=15; A="Make sure to lock the door."; B="문 꼭 잠그도록 해."},
    @{Row=16; A="I'm looking forward to the party."; B="파티가 정말 기대돼."},
    @{Row=17; A="Can I get a glass of water?"; B="물 한 잔 주시겠어요?"},
    @{Row=18; A="It takes about an hour to get there."; B="거기까지 가는 데 한 시간 정도 걸려."},
    @{Row=19; A="That's why I was late."; B="그래서 늦은 거야."},
    @{Row=20; A="I have no idea what you're talking about."; B="무슨 말을 하는지 전혀 모르겠어."},
    @{Row=21; A="Let me check my schedule."; B="내 일정 좀 확인해볼게."},
    @{Row=22; A="I feel like having pizza today."; B="오늘 피자 먹고 싶은 기분이야."},
    @{Row=23; A="I'm afraid I can't help you."; B="유감스럽지만 널 도와줄 수 없을 것 같아."},
    @{Row=24; A="You'd better see a doctor."; B="병원에 가보는 게 좋을 거야."},
    @{Row=25; A="It's hard to believe."; B="믿기 힘들어."},
    @{Row=26; A="Is it okay to use this phone?"; B="이 전화기 써도 괜찮아?"},
    @{Row=27; A="Thank you for inviting me."; B="초대해줘서 고마워."},
    @{Row=28; A="I'm thinking of moving to Seoul."; B="서울로 이사 갈까 생각 중이야."},
    @{Row=29; A="What makes you think so?"; B="왜 그렇게 생각하는 거야?"}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
}
